$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige les noms des paroisses (remplace le trait d'union par un tiret demi-cadratin
# entoure d'espaces). L'ordre d'ecriture (lignes 2, 3 puis 5) reproduit l'ordre dans
# lequel les chaines partagees apparaissent dans le classeur final.
$ws.Range("E2").Value2 = "Pully – Paudex"
$ws.Range("E3").Value2 = "Belmont – Lutry"
$ws.Range("E5").Value2 = "Savigny – Forel"

# Met a jour la selection active de la feuille.
$ws.Range("E5").Select() | Out-Null
